$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("workflow")

# Update fourier terms value (B3): 2 -> 0
$ws.Range("B3").Value = 0

# Add the new "random slopes" row (row 7) before changing B4 so that the
# shared-string table fills in the same order as the authoring commit:
#   53 random slopes, 54 list_rand_slopes, 55 (1|store), 56 (TV1|store), (TV2|store)
$ws.Range("A7").Value = "random slopes"
$ws.Range("C7").Value = "list_rand_slopes"

# Update the random effects intercepts example value (B4): store -> (1|store)
$ws.Range("B4").Value = "(1|store)"

# Fill in the new random slopes example value (B7)
$ws.Range("B7").Value = "(TV1|store), (TV2|store)"

# Make "workflow" the active/selected sheet (was "variables"), with B4 selected
$ws.Activate()
$ws.Range("B4").Select() | Out-Null
